$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.074.34"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "'2.917.83"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'594.93"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("D6").Value = "'143.32"
$ws.Range("E6").Value = "  -1.60%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "'0.499"
$ws.Range("E8").Value = "  -1.48%  "

$ws.Range("D9").Value = "'6.91"
$ws.Range("E9").Value = "  +1.49%  "

$ws.Range("E10").Value = "  -2.22%  "

$ws.Range("D11").Value = "'0.435"
$ws.Range("E11").Value = "  -0.91%  "

$ws.Range("D12").Value = "'0.0000222"
$ws.Range("E12").Value = "  -1.28%  "

$ws.Range("D13").Value = "'33.04"
$ws.Range("E13").Value = "  -1.72%  "

$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").Value = "'3.407.93"
$ws.Range("E15").Value = "  +0.16%  "

$ws.Range("D16").Value = "'61.103.69"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").Value = "'2.919.93"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "'6.62"
$ws.Range("E18").Value = "  -1.21%  "

$ws.Range("D19").Value = "'430.95"
$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").Value = "'13.47"
$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").Value = "'0.670"
$ws.Range("E21").Value = "  -1.68%  "

$ws.Range("D22").Value = "'7.01"
$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("D23").Value = "'81.60"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'10.88"
$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("D25").Value = "'2.16"
$ws.Range("E25").Value = "  -2.52%  "

$ws.Range("D26").Value = "'11.67"
$ws.Range("E26").Value = "  -2.10%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "'2.19"
$ws.Range("E28").Value = "  -4.80%  "

$ws.Range("D29").Value = "'2.58"
$ws.Range("E29").Value = "  -1.10%  "

$ws.Range("D30").Value = "'6.86"
$ws.Range("E30").Value = "  -2.64%  "

$ws.Range("D31").Value = "'26.48"
$ws.Range("E31").Value = "  +0.37%  "

$ws.Range("D32").Value = "'0.108"
$ws.Range("E32").Value = "  +1.03%  "

$ws.Range("E33").Value = "  +0.17%  "

$ws.Range("D34").Value = "'0.0" + [char]8323 + "0869"
$ws.Range("E34").Value = "  +2.32%  "

$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("D36").Value = "'5.58"
$ws.Range("E36").Value = "  -0.50%  "

$ws.Range("D37").Value = "'2.94"
$ws.Range("E37").Value = "  -2.45%  "

$ws.Range("D38").Value = "'1.98"
$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("E39").Value = "  -0.50%  "

$ws.Range("D40").Value = "'8.49"
$ws.Range("E40").Value = "  -0.79%  "

$ws.Range("D41").Value = "'42.08"
$ws.Range("E41").Value = "  +4.38%  "

$ws.Range("D42").Value = "'0.278"
$ws.Range("E42").Value = "  -2.87%  "

$ws.Range("D43").Value = "'0.0343"
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("D44").Value = "'2.683.18"
$ws.Range("E44").Value = "  -0.52%  "

$ws.Range("D45").Value = "'133.07"
$ws.Range("E45").Value = "  +1.65%  "

$ws.Range("D46").Value = "'361.59"
$ws.Range("E46").Value = "  -3.25%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Value = "'23.35"
$ws.Range("E48").Value = "  -2.43%  "

$ws.Range("D49").Value = "'0.104"
$ws.Range("E49").Value = "  -1.55%  "

$ws.Range("D50").Value = "'1.98"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("D51").Value = "'0.124"
$ws.Range("E51").Value = "  -0.78%  "
